$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Dmitri"
$ws.Range("B2").Value = "Dmitri Nikolaenko"
$ws.Range("C2").Value = "Njord Law"
$ws.Range("D2").Value = "Latvia"
$ws.Range("G2").Value = "dn@njordlaw.lv"

$ws.Range("A5").Value = "Giuseppe"
$ws.Range("B5").Value = "Giuseppe Abbruzzese"
$ws.Range("C5").Value = "Legance"
$ws.Range("D5").Value = "Italy"
$ws.Range("G5").Value = "gabbruzzese@legance.it"

$ws.Range("A6").Value = "Shiran"
$ws.Range("B6").Value = "Shiran Sofer"
$ws.Range("C6").Value = "Gornitzky And Co"
$ws.Range("D6").Value = "Israel"
$ws.Range("G6").Value = "shirans@gornitzky.com"

$ws.Range("A7").Value = "Robert"
$ws.Range("B7").Value = "Robert Nader"
$ws.Range("C7").Value = "Forbes Hare"
$ws.Range("D7").Value = "British Virgin Islands"
$ws.Range("G7").Value = "robert.nader@forbeshare.com"

$ws.Range("A10").Value = "Adi"
$ws.Range("B10").Value = "Adi Ron"
$ws.Range("C10").Value = "Fischer"
$ws.Range("D10").Value = "Israel"
$ws.Range("G10").Value = "aron@fbclawyers.com"

$ws.Range("A11").Value = "Qiuzhao"
$ws.Range("B11").Value = "Qiuzhao Wang"
$ws.Range("C11").Value = "TC Law Firm"
$ws.Range("D11").Value = "Hong Kong"
$ws.Range("G11").Value = "wqc@tclawfirm.com"

$ws.Range("A13").Value = "Malek"
$ws.Range("B13").Value = "Malek Barudi Mjur (Oxford)"
$ws.Range("C13").Value = "Taylor Wessing"
$ws.Range("D13").Value = "Germany"
$ws.Range("G13").Value = "m.barudi@taylorwessing.com"

$ws.Range("A15").Value = "David"
$ws.Range("B15").Value = "David Goldman"
$ws.Range("C15").Value = "Gornitzky And Co"
$ws.Range("D15").Value = "Israel"
$ws.Range("G15").Value = "davidg@gornitzky.com"

$ws.Range("A17").Value = "Giuseppe"
$ws.Range("B17").Value = "Giuseppe Abbruzzese"
$ws.Range("C17").Value = "Legance"
$ws.Range("D17").Value = "Italy"
$ws.Range("G17").Value = "gabbruzzese@legance.it"

$ws.Range("A18").Value = "Stefanos"
$ws.Range("B18").Value = "Stefanos Charaktiniotis"
$ws.Range("C18").Value = "Zepos And Yannopoulos"
$ws.Range("D18").Value = "Greece"
$ws.Range("G18").Value = "s.charaktiniotis@zeya.com"

$ws.Range("A19").Value = "Martin"
$ws.Range("B19").Value = "Martin Leboutillier"
$ws.Range("C19").Value = "Collas Crill"
$ws.Range("D19").Value = "Bermuda"
$ws.Range("G19").Value = "martin.leboutillier@collascrill.com"

$ws.Range("A21").Value = "Andrea"
$ws.Range("B21").Value = "Andrea Calvi"
$ws.Range("C21").Value = "Pedersoli"
$ws.Range("D21").Value = "Italy"
$ws.Range("G21").Value = "acalvi@pglex.it"

$ws.Range("A23").Value = "Kevin"
$ws.Range("B23").Value = "Kevin Tsen"
$ws.Range("C23").Value = "CFN Law"
$ws.Range("D23").Value = "Hong Kong"
$ws.Range("G23").Value = "kevin.tsen@cfnlaw.com.hk"

$ws.Range("A25").Value = "Simone"
$ws.Range("B25").Value = "Simone Ambrogi"
$ws.Range("C25").Value = "Legance"
$ws.Range("D25").Value = "Italy"
$ws.Range("G25").Value = "sambrogi@legance.it"

$ws.Range("A26").Value = "Ya-Chiao"
$ws.Range("B26").Value = "Ya-Chiao Chang"
$ws.Range("C26").Value = "Winston And Strawn"
$ws.Range("D26").Value = "China"
$ws.Range("G26").Value = "ychang@winston.com"

$ws.Range("A27").Value = "Andrew"
$ws.Range("B27").Value = "Andrew Feighery"
$ws.Range("C27").Value = "Byrne Wallace"
$ws.Range("D27").Value = "Ireland"
$ws.Range("G27").Value = "afeighery@byrnewallace.com"

$ws.Range("A28").Value = "Catherine"
$ws.Range("B28").Value = "Catherine Ross"
$ws.Range("C28").Value = "Forbes Hare"
$ws.Range("D28").Value = "Singapore"
$ws.Range("G28").Value = "catherine.ross@forbeshare.com"

$ws.Range("A29").Value = "Hartwig"
$ws.Range("B29").Value = "Hartwig Kienast"
$ws.Range("C29").Value = "Wolf Theiss"
$ws.Range("D29").Value = "Austria"
$ws.Range("G29").Value = "hartwig.kienast@wolftheiss.com"

$ws.Range("A31").Value = "Henning"
$ws.Range("B31").Value = "Henning Von Lillienskjold"
$ws.Range("C31").Value = "DahlLaw"
$ws.Range("D31").Value = "Denmark"
$ws.Range("G31").Value = "hvl@dahllaw.dk"

$ws.Range("A32").Value = "Achiron"
$ws.Range("B32").Value = "Achiron Jonathan"
$ws.Range("C32").Value = "EBN"
$ws.Range("D32").Value = "Israel"
$ws.Range("G32").Value = "jonathana@ebnlaw.co.il"

$ws.Range("A33").Value = "Mark"
$ws.Range("B33").Value = "Mark Davis"
$ws.Range("C33").Value = "Mishcon Karas"
$ws.Range("D33").Value = "Hong Kong"
$ws.Range("G33").Value = "mark.davis@mishcon.com"

$ws.Range("A34").Value = "Wang"
$ws.Range("B34").Value = "Wang Haiyu (Nicole Wang)"
$ws.Range("C34").Value = "Longan Law"
$ws.Range("D34").Value = "China"
$ws.Range("G34").Value = "wanghaiyu@longanlaw.com"

$ws.Range("A36").Value = "Bai"
$ws.Range("B36").Value = "Bai Xianyue"
$ws.Range("C36").Value = "Grandall"
$ws.Range("D36").Value = "China"
$ws.Range("G36").Value = "baixianyue@grandall.com.cn"

$ws.Range("A37").Value = "Christopher"
$ws.Range("B37").Value = "Christopher Bromilow"
$ws.Range("C37").Value = "Forbes Hare"
$ws.Range("D37").Value = "British Virgin Islands"
$ws.Range("G37").Value = "christopher.bromilow@forbeshare.com"

$ws.Range("A38").Value = "Jingzhong"
$ws.Range("B38").Value = "Jingzhong Zhang"
$ws.Range("C38").Value = "TC Law Firm"
$ws.Range("D38").Value = "China"
$ws.Range("G38").Value = "undefined zhangjz@tclawfirm.com"

$ws.Range("A42").Value = "Barry"
$ws.Range("B42").Value = "Barry Smith"
$ws.Range("C42").Value = "Forbes Hare"
$ws.Range("D42").Value = "Cayman Islands"
$ws.Range("G42").Value = "barry.smith@forbeshare.com"

$ws.Range("A43").Value = "Abramovich"
$ws.Range("B43").Value = "Abramovich Menachem"
$ws.Range("C43").Value = "EBN"
$ws.Range("D43").Value = "Israel"
$ws.Range("G43").Value = "menachema@ebnlaw.co.il"

$ws.Range("A44").Value = "William"
$ws.Range("B44").Value = "William Hare"
$ws.Range("C44").Value = "Forbes Hare"
$ws.Range("D44").Value = "British Virgin Islands"
$ws.Range("G44").Value = "william.hare@forbeshare.com"

$ws.Range("A45").Value = "Feng"
$ws.Range("B45").Value = "Feng Tao"
$ws.Range("C45").Value = "Longan Law"
$ws.Range("D45").Value = "China"
$ws.Range("G45").Value = "fengtao@longanlaw.com"

$ws.Range("A47").Value = "Adam"
$ws.Range("B47").Value = "Adam Kadesh"
$ws.Range("C47").Value = "Fischer"
$ws.Range("D47").Value = "Israel"
$ws.Range("G47").Value = "akadesh@fbclawyers.com"

$ws.Range("A49").Value = "Fang"
$ws.Range("B49").Value = "Fang Ling"
$ws.Range("C49").Value = "Longan Law"
$ws.Range("D49").Value = "China"
$ws.Range("G49").Value = "fangling@longanlaw.com"

$ws.Range("A50").Value = "Nitzan"
$ws.Range("B50").Value = "Nitzan Aberbach"
$ws.Range("C50").Value = "EBN"
$ws.Range("D50").Value = "Israel"
$ws.Range("G50").Value = "nitzana@ebnlaw.co.il"

$ws.Range("A52").Value = "Ido"
$ws.Range("B52").Value = "Ido Malin"
$ws.Range("C52").Value = "Gornitzky And Co"
$ws.Range("D52").Value = "Israel"
$ws.Range("G52").Value = "idom@gornitzky.com"

